# Add a new "looping" entry (row 17) with value 5 under the existing
# account id list, then move the selection down to the next empty cell
# (A18), matching where the user would continue entering data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 5

$ws.Range("A18").Select() | Out-Null
